# Apply updated odds values to Sheet1 per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of row -> hashtable of column letter -> new value
$changes = @{
    2  = @{ "G" = 2.55; "I" = 2.7; "AJ" = 29 }
    3  = @{ "J" = 1.03; "K" = 17; "L" = 1.14; "M" = 5.5; "V" = 9.5; "W" = 12 }
    15 = @{
        "G" = 1.5; "H" = 3.65; "I" = 5.8; "N" = 1.8; "O" = 1.8;
        "U" = 5.8; "V" = 6.8; "X" = 10.25; "Y" = 21; "AA" = 6.3;
        "AB" = 14; "AC" = 60; "AD" = 450; "AE" = 12.5; "AF" = 29;
        "AH" = 90; "AI" = 50
    }
    16 = @{
        "G" = 1.95; "H" = 3.05; "I" = 3.7; "N" = 2.12;
        "U" = 7.3; "V" = 7.3; "W" = 14; "X" = 14; "Y" = 25;
        "Z" = 7.4; "AB" = 13; "AE" = 7.6; "AF" = 15; "AG" = 10.75;
        "AH" = 45; "AI" = 30; "AJ" = 37
    }
    19 = @{
        "G" = 1.09; "H" = 7.3; "I" = 17;
        "U" = 5.5; "V" = 10.5; "W" = 5.2; "X" = 9.75;
        "Z" = 19; "AA" = 15.5; "AB" = 35; "AF" = 150; "AI" = 250; "AJ" = 175
    }
    30 = @{
        "G" = 1.24; "H" = 5.6; "J" = 1.02; "K" = 10; "L" = 1.12;
        "M" = 5.3; "N" = 1.4; "O" = 2.75; "P" = 1.23; "Q" = 3.75;
        "R" = 1.8; "S" = 1.91; "T" = 9.75; "U" = 7.4;
        "Z" = 10; "AA" = 12; "AB" = 20; "AC" = 75; "AD" = 450;
        "AE" = 37; "AF" = 80; "AI" = 110; "AJ" = 75
    }
}

foreach ($rowNum in $changes.Keys) {
    $cols = $changes[$rowNum]
    foreach ($col in $cols.Keys) {
        $addr = "$col$rowNum"
        $ws.Range($addr).Value = $cols[$col]
    }
}
